# Insert a new weekly data row at row 59 (Fruta / hortaliza, semanal).
# Excel shifts the existing rows 59-157 down to 60-158 and we populate the
# freshly-inserted row 59 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("59").Insert()

$ws.Range("A59").Value = 7
$ws.Range("B59").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C59").Value = "Ñuble"
$ws.Range("D59").Value = 45117
$ws.Range("E59").Value = 16
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100108
$ws.Range("H59").Value = "Tropicales y subtropicales"
$ws.Range("I59").Value = 100108002
$ws.Range("J59").Value = "Mango"
$ws.Range("K59").Value = "Sin especificar"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 60
$ws.Range("N59").Value = 9000
$ws.Range("O59").Value = 9000
$ws.Range("P59").Value = 9000
$ws.Range("Q59").Value = "$/bandeja 4 kilos"
$ws.Range("R59").Value = "Brasil"
$ws.Range("S59").Value = 2250
$ws.Range("T59").Value = 4
